$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B2").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Test"
$ws.Range("B5").Value = "Test2"
$ws.Range("C5").Value = 34.018956305811898
$ws.Range("D5").Value = -118.28375294545

$ws.Range("D11").Select()
